# Update GSC export data:
#  - Append 4 new daily rows (2025-11-05 .. 2025-11-08) to the "Chart" sheet
#  - Update two values on the "Critical issues" sheet

$wb = $excel.ActiveWorkbook

# --- Chart sheet: append rows 33-36 -----------------------------------
$chart = $wb.Worksheets.Item("Chart")

$newRows = @(
    @("2025-11-05", 102.0, 205.0, 31.0),
    @("2025-11-06", 102.0, 205.0, 34.0),
    @("2025-11-07", 102.0, 205.0, 23.0),
    @("2025-11-08", 102.0, 205.0, 17.0)
)

$startRow = 33
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $chart.Cells.Item($r, 1).Value = $row[0]
    $chart.Cells.Item($r, 2).Value = $row[1]
    $chart.Cells.Item($r, 3).Value = $row[2]
    $chart.Cells.Item($r, 4).Value = $row[3]
}

# --- Critical issues sheet: update two data points --------------------
$critical = $wb.Worksheets.Item("Critical issues")
$critical.Cells.Item(2, 4).Value = 26.0
$critical.Cells.Item(6, 4).Value = 41.0
